# --- Step 1: Update "总计" (summary) sheet with new 2022-Q3 row ---
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Copy row5 formatting down to new row6 (keeps A6 styled like A5, s=2)
$ws1.Range("A5:D5").Copy($ws1.Range("A6:D6"))

# Shift the B/C/D content down by one row (index column A stays fixed per-row),
# working bottom-up so source cells aren't overwritten before being read
$ws1.Cells.Item(6,2).Value = $ws1.Cells.Item(5,2).Value2
$ws1.Cells.Item(6,3).Value = $ws1.Cells.Item(5,3).Value2
$ws1.Cells.Item(6,4).Value = $ws1.Cells.Item(5,4).Value2
$ws1.Cells.Item(6,1).Value = 4

$ws1.Cells.Item(5,2).Value = $ws1.Cells.Item(4,2).Value2
$ws1.Cells.Item(5,3).Value = $ws1.Cells.Item(4,3).Value2
$ws1.Cells.Item(5,4).Value = $ws1.Cells.Item(4,4).Value2

$ws1.Cells.Item(4,2).Value = $ws1.Cells.Item(3,2).Value2
$ws1.Cells.Item(4,3).Value = $ws1.Cells.Item(3,3).Value2
$ws1.Cells.Item(4,4).Value = $ws1.Cells.Item(3,4).Value2

$ws1.Cells.Item(3,2).Value = $ws1.Cells.Item(2,2).Value2
$ws1.Cells.Item(3,3).Value = $ws1.Cells.Item(2,3).Value2
$ws1.Cells.Item(3,4).Value = $ws1.Cells.Item(2,4).Value2

$ws1.Cells.Item(2,2).Value = "2022-Q3"
$ws1.Cells.Item(2,3).Value = 8
$ws1.Cells.Item(2,4).Value = 12.83

# --- Step 2: Add new "2022-Q3" worksheet (copy of "2022-Q2" sheet, before it) ---
$target = $wb.Worksheets.Item(2)
$target.Copy($target)
$ws = $wb.Worksheets.Item(2)
$ws.Name = "2022-Q3"
$ws.Range("A8:H8").Copy($ws.Range("A9:H9"))

$ws.Cells.Item(1,2).Value = "基金代码"
$ws.Cells.Item(1,3).Value = "基金名称"
$ws.Cells.Item(1,4).Value = "基金规模"
$ws.Cells.Item(1,5).Value = "股票总仓位"
$ws.Cells.Item(1,6).Value = "仓位占比"
$ws.Cells.Item(1,7).Value = "持有市值(亿元)"
$ws.Cells.Item(1,8).Value = "仓位排名"
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "'008903"
$ws.Cells.Item(2,3).Value = "广发科技先锋混合"
$ws.Cells.Item(2,4).Value = "'134.64"
$ws.Cells.Item(2,5).Value = "'94.69"
$ws.Cells.Item(2,6).Value = "'4.44"
$ws.Cells.Item(2,7).Value = "'5.9780"
$ws.Cells.Item(2,8).Value = 9
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "'162703"
$ws.Cells.Item(3,3).Value = "广发小盘成长混合（LOF）A"
$ws.Cells.Item(3,4).Value = "'99.59"
$ws.Cells.Item(3,5).Value = "'88.79"
$ws.Cells.Item(3,6).Value = "'4.69"
$ws.Cells.Item(3,7).Value = "'4.6708"
$ws.Cells.Item(3,8).Value = 9
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "'003745"
$ws.Cells.Item(4,3).Value = "广发多元新兴股票"
$ws.Cells.Item(4,4).Value = "'35.14"
$ws.Cells.Item(4,5).Value = "'90.17"
$ws.Cells.Item(4,6).Value = "'4.30"
$ws.Cells.Item(4,7).Value = "'1.5110"
$ws.Cells.Item(4,8).Value = 9
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "'009086"
$ws.Cells.Item(5,3).Value = "鹏华价值共赢两年持有期混合"
$ws.Cells.Item(5,4).Value = "'11.02"
$ws.Cells.Item(5,5).Value = "'89.84"
$ws.Cells.Item(5,6).Value = "'3.30"
$ws.Cells.Item(5,7).Value = "'0.3637"
$ws.Cells.Item(5,8).Value = 10
$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "'009132"
$ws.Cells.Item(6,3).Value = "广发小盘成长混合（LOF）C"
$ws.Cells.Item(6,4).Value = "'5.90"
$ws.Cells.Item(6,5).Value = "'88.79"
$ws.Cells.Item(6,6).Value = "'4.69"
$ws.Cells.Item(6,7).Value = "'0.2767"
$ws.Cells.Item(6,8).Value = 9
$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "'012272"
$ws.Cells.Item(7,3).Value = "渤海汇金创新价值一年持有期混合型发起式证券投资基金"
$ws.Cells.Item(7,4).Value = "'0.98"
$ws.Cells.Item(7,5).Value = "'87.72"
$ws.Cells.Item(7,6).Value = "'2.86"
$ws.Cells.Item(7,7).Value = "'0.0280"
$ws.Cells.Item(7,8).Value = 6
$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "'015225"
$ws.Cells.Item(8,3).Value = "汇添富中证细分化工产业主题指数增强A"
$ws.Cells.Item(8,4).Value = "'0.11"
$ws.Cells.Item(8,5).Value = "'92.25"
$ws.Cells.Item(8,6).Value = "'2.67"
$ws.Cells.Item(8,7).Value = "'0.0029"
$ws.Cells.Item(8,8).Value = 9
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "'015226"
$ws.Cells.Item(9,3).Value = "汇添富中证细分化工产业主题指数增强C"
$ws.Cells.Item(9,4).Value = "'0.02"
$ws.Cells.Item(9,5).Value = "'92.25"
$ws.Cells.Item(9,6).Value = "'2.67"
$ws.Cells.Item(9,7).Value = "'0.0005"
$ws.Cells.Item(9,8).Value = 9

Write-Host "Edit complete."
